# Actualizacion del cronograma y Hoja de desempeño
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Cronograma: marcar como completadas (100%) las tareas de las filas 16 y 24
$ws.Range("H16").Value = 1
$ws.Range("H24").Value = 1

# Hoja de desempeño: el promedio general (H25) se recalcula solo a partir
# de los valores anteriores; dejamos la selección activa donde quedó el
# usuario tras revisar el resultado.
[void]$ws.Range("H25").Select()
